$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing cells per diff ---
$ws.Cells.Item(56, 17).Value = 0   # Q56: 1 -> 0
$ws.Cells.Item(60, 17).Value = 0   # Q60: 2 -> 0
$ws.Cells.Item(69, 17).Value = 0   # Q69: 2 -> 0
$ws.Cells.Item(71, 17).Value = 0   # Q71: 1 -> 0

$ws.Cells.Item(1485, 15).Value = 2 # O1485: 0 -> 2
$ws.Cells.Item(1487, 15).Value = 1 # O1487: 0 -> 1
$ws.Cells.Item(1487, 18).Value = 0 # R1487: blank -> 0
$ws.Cells.Item(1488, 18).Value = 0 # R1488: blank -> 0

# --- Append new weekly rows 1489-1506 ---
$row = 1489
$ws.Cells.Item($row, 1).Value = 45474
$ws.Cells.Item($row, 2).Value = 3382
$ws.Cells.Item($row, 3).Value = 3450
$ws.Cells.Item($row, 4).Value = 3257
$ws.Cells.Item($row, 5).Value = 3269.35009765625
$ws.Cells.Item($row, 6).Value = 3269.35009765625
$ws.Cells.Item($row, 7).Value = 5603909
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 7
$ws.Cells.Item($row, 10).Value = 1
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 27
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1490
$ws.Cells.Item($row, 1).Value = 45481
$ws.Cells.Item($row, 2).Value = 3204
$ws.Cells.Item($row, 3).Value = 3269
$ws.Cells.Item($row, 4).Value = 3126.10009765625
$ws.Cells.Item($row, 5).Value = 3229.85009765625
$ws.Cells.Item($row, 6).Value = 3229.85009765625
$ws.Cells.Item($row, 7).Value = 11407642
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 7
$ws.Cells.Item($row, 10).Value = 8
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 28
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1491
$ws.Cells.Item($row, 1).Value = 45488
$ws.Cells.Item($row, 2).Value = 3241.75
$ws.Cells.Item($row, 3).Value = 3279
$ws.Cells.Item($row, 4).Value = 3205.39990234375
$ws.Cells.Item($row, 5).Value = 3259
$ws.Cells.Item($row, 6).Value = 3259
$ws.Cells.Item($row, 7).Value = 3942711
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 7
$ws.Cells.Item($row, 10).Value = 15
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 29
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 1
$ws.Cells.Item($row, 17).Value = 1
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1492
$ws.Cells.Item($row, 1).Value = 45495
$ws.Cells.Item($row, 2).Value = 3250
$ws.Cells.Item($row, 3).Value = 3552.5
$ws.Cells.Item($row, 4).Value = 3223.199951171875
$ws.Cells.Item($row, 5).Value = 3494.14990234375
$ws.Cells.Item($row, 6).Value = 3494.14990234375
$ws.Cells.Item($row, 7).Value = 15025069
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 7
$ws.Cells.Item($row, 10).Value = 22
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 30
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1493
$ws.Cells.Item($row, 1).Value = 45502
$ws.Cells.Item($row, 2).Value = 3495
$ws.Cells.Item($row, 3).Value = 3499.89990234375
$ws.Cells.Item($row, 4).Value = 3388
$ws.Cells.Item($row, 5).Value = 3462.35009765625
$ws.Cells.Item($row, 6).Value = 3462.35009765625
$ws.Cells.Item($row, 7).Value = 6939362
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 7
$ws.Cells.Item($row, 10).Value = 29
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 31
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 2
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1494
$ws.Cells.Item($row, 1).Value = 45509
$ws.Cells.Item($row, 2).Value = 3320.050048828125
$ws.Cells.Item($row, 3).Value = 3459
$ws.Cells.Item($row, 4).Value = 3283.89990234375
$ws.Cells.Item($row, 5).Value = 3331.699951171875
$ws.Cells.Item($row, 6).Value = 3331.699951171875
$ws.Cells.Item($row, 7).Value = 6486280
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 8
$ws.Cells.Item($row, 10).Value = 5
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 32
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1495
$ws.Cells.Item($row, 1).Value = 45516
$ws.Cells.Item($row, 2).Value = 3331.64990234375
$ws.Cells.Item($row, 3).Value = 3450.10009765625
$ws.Cells.Item($row, 4).Value = 3295.35009765625
$ws.Cells.Item($row, 5).Value = 3444.75
$ws.Cells.Item($row, 6).Value = 3444.75
$ws.Cells.Item($row, 7).Value = 4348754
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 8
$ws.Cells.Item($row, 10).Value = 12
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 33
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1496
$ws.Cells.Item($row, 1).Value = 45523
$ws.Cells.Item($row, 2).Value = 3436.199951171875
$ws.Cells.Item($row, 3).Value = 3625
$ws.Cells.Item($row, 4).Value = 3436.199951171875
$ws.Cells.Item($row, 5).Value = 3570
$ws.Cells.Item($row, 6).Value = 3570
$ws.Cells.Item($row, 7).Value = 6356761
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 8
$ws.Cells.Item($row, 10).Value = 19
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 34
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1497
$ws.Cells.Item($row, 1).Value = 45530
$ws.Cells.Item($row, 2).Value = 3592
$ws.Cells.Item($row, 3).Value = 3658.800048828125
$ws.Cells.Item($row, 4).Value = 3510
$ws.Cells.Item($row, 5).Value = 3565.14990234375
$ws.Cells.Item($row, 6).Value = 3565.14990234375
$ws.Cells.Item($row, 7).Value = 6097535
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 8
$ws.Cells.Item($row, 10).Value = 26
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 35
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1498
$ws.Cells.Item($row, 1).Value = 45537
$ws.Cells.Item($row, 2).Value = 3565.14990234375
$ws.Cells.Item($row, 3).Value = 3753.949951171875
$ws.Cells.Item($row, 4).Value = 3552
$ws.Cells.Item($row, 5).Value = 3695.5
$ws.Cells.Item($row, 6).Value = 3695.5
$ws.Cells.Item($row, 7).Value = 7207950
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 9
$ws.Cells.Item($row, 10).Value = 2
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 36
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1499
$ws.Cells.Item($row, 1).Value = 45544
$ws.Cells.Item($row, 2).Value = 3690.10009765625
$ws.Cells.Item($row, 3).Value = 3799.85009765625
$ws.Cells.Item($row, 4).Value = 3653.75
$ws.Cells.Item($row, 5).Value = 3767
$ws.Cells.Item($row, 6).Value = 3767
$ws.Cells.Item($row, 7).Value = 5204927
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 9
$ws.Cells.Item($row, 10).Value = 9
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 37
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1500
$ws.Cells.Item($row, 1).Value = 45551
$ws.Cells.Item($row, 2).Value = 3768
$ws.Cells.Item($row, 3).Value = 3819
$ws.Cells.Item($row, 4).Value = 3708
$ws.Cells.Item($row, 5).Value = 3797.199951171875
$ws.Cells.Item($row, 6).Value = 3797.199951171875
$ws.Cells.Item($row, 7).Value = 4096591
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 9
$ws.Cells.Item($row, 10).Value = 16
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 38
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1501
$ws.Cells.Item($row, 1).Value = 45558
$ws.Cells.Item($row, 2).Value = 3797.199951171875
$ws.Cells.Item($row, 3).Value = 3867
$ws.Cells.Item($row, 4).Value = 3701.39990234375
$ws.Cells.Item($row, 5).Value = 3816.699951171875
$ws.Cells.Item($row, 6).Value = 3816.699951171875
$ws.Cells.Item($row, 7).Value = 6782611
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 9
$ws.Cells.Item($row, 10).Value = 23
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 39
$ws.Cells.Item($row, 15).Value = 1
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1502
$ws.Cells.Item($row, 1).Value = 45565
$ws.Cells.Item($row, 2).Value = 3786
$ws.Cells.Item($row, 3).Value = 3863.550048828125
$ws.Cells.Item($row, 4).Value = 3625
$ws.Cells.Item($row, 5).Value = 3670.10009765625
$ws.Cells.Item($row, 6).Value = 3670.10009765625
$ws.Cells.Item($row, 7).Value = 4983377
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 9
$ws.Cells.Item($row, 10).Value = 30
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 40
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1503
$ws.Cells.Item($row, 1).Value = 45572
$ws.Cells.Item($row, 2).Value = 3739.949951171875
$ws.Cells.Item($row, 3).Value = 3748
$ws.Cells.Item($row, 4).Value = 3415.10009765625
$ws.Cells.Item($row, 5).Value = 3474.39990234375
$ws.Cells.Item($row, 6).Value = 3474.39990234375
$ws.Cells.Item($row, 7).Value = 7294034
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 10
$ws.Cells.Item($row, 10).Value = 7
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 41
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1504
$ws.Cells.Item($row, 1).Value = 45579
$ws.Cells.Item($row, 2).Value = 3474.39990234375
$ws.Cells.Item($row, 3).Value = 3529.449951171875
$ws.Cells.Item($row, 4).Value = 3296.10009765625
$ws.Cells.Item($row, 5).Value = 3381.449951171875
$ws.Cells.Item($row, 6).Value = 3381.449951171875
$ws.Cells.Item($row, 7).Value = 3573950
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 10
$ws.Cells.Item($row, 10).Value = 14
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 42
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1505
$ws.Cells.Item($row, 1).Value = 45586
$ws.Cells.Item($row, 2).Value = 3417.949951171875
$ws.Cells.Item($row, 3).Value = 3421.5
$ws.Cells.Item($row, 4).Value = 3232.35009765625
$ws.Cells.Item($row, 5).Value = 3266.550048828125
$ws.Cells.Item($row, 6).Value = 3266.550048828125
$ws.Cells.Item($row, 7).Value = 4044051
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 10
$ws.Cells.Item($row, 10).Value = 21
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 43
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row = 1506
$ws.Cells.Item($row, 1).Value = 45593
$ws.Cells.Item($row, 2).Value = 3270
$ws.Cells.Item($row, 3).Value = 3321.800048828125
$ws.Cells.Item($row, 4).Value = 3227
$ws.Cells.Item($row, 5).Value = 3267.050048828125
$ws.Cells.Item($row, 6).Value = 3267.050048828125
$ws.Cells.Item($row, 7).Value = 3503250
$ws.Cells.Item($row, 8).Value = 2024
$ws.Cells.Item($row, 9).Value = 10
$ws.Cells.Item($row, 10).Value = 28
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 44
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0
$ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

